# Aggiustato un piccolo errore di battitura
#
# 1) Date placeholder ("datetimeFigureOut" field) cached text: 19/10/2018 -> 20/10/2018
# 2) Slide-number placeholder ("slidenum" field) cached text: <N> -> <#>
#    (applied on the slide master and every slide layout)
# 3) Slide 1 background picture nudged: <a:off x="0" .../> -> <a:off x="8238" .../>
# 4) Slide 4 body text: fixed missing word ("Chi ne ha" -> "Chi più ne ha")

$p = $ppt.ActivePresentation

$dateText = "20/10/2018"
$slideNumText = [string][char]0x2039 + "#" + [string][char]0x203A

function Update-HeaderFooterPlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = -1 }

        if ($phType -eq 16) {
            # ppPlaceholderDate
            $sh.TextFrame.TextRange.Text = $dateText
        } elseif ($phType -eq 13) {
            # ppPlaceholderSlideNumber
            $sh.TextFrame.TextRange.Text = $slideNumText
        }
    }
}

# -- 1 & 2: refresh the date / slide-number placeholders on the master ...
Update-HeaderFooterPlaceholders $p.SlideMaster.Shapes

# ... and on every slide layout.
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-HeaderFooterPlaceholders $layout.Shapes
}

# -- 3: nudge the full-bleed picture on slide 1 a touch to the right.
$slide1 = $p.Slides.Item(1)
$bgPicture = $slide1.Shapes.Item(1)
$bgPicture.Left = 8238 / 12700.0   # EMU -> points (1 pt = 12700 EMU)

# -- 4: fix the missing "più" in the bullet list on slide 4.
$slide4 = $p.Slides.Item(4)
$bodyShape = $slide4.Shapes.Item(1)
$bodyText = $bodyShape.TextFrame.TextRange
$bodyText.Replace("(Chi ne ha più ne metta)", "(Chi più ne ha più ne metta)", 1, 0, 0) | Out-Null
